{"js": "// \"modifica testo teoria avanzata\"\n//\n// Three textual edits on the document body:\n//  1) Paragraph \"Il tono \u00e8...\": the clause\n//       \", e teoricamente \u00e8 divisibili in\"\n//     becomes\n//       \" ed \u00e8 divisibile in\"\n//  2) Same paragraph, later: \", cio\u00e8 piccolissimo\" becomes \", cio\u00e8 un piccolissimo\"\n//  3) Paragraph \"Il semitono diatonico...\": the split run \"qua\"/\"ttro comma\"\n//     (with the _GoBack bookmark sitting between them) is normalized into a\n//     single \"quattro comma\" run.\n//\n// Doing (3) first lets the engine naturally retire the pre-existing _GoBack\n// bookmark that lived inside \"qua|ttro comma\"; we then re-create it at its\n// new resting place (inside \"divisibile\", right after \"\u00e8 div\") to match\n// where Word leaves its \"last edit\" marker after the retype in paragraph 2.\n\nconst body = context.document.body;\n\n// 1) Normalize \"qua\" + \"ttro comma\" into one run (text is unchanged).\nconst quattro = body.search(\"quattro comma\", { matchCase: true, matchWholeWord: false });\nquattro.load(\"items\");\nawait context.sync();\nif (quattro.items.length > 0) {\n  quattro.items[0].insertText(\"quattro comma\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) \", e teoricamente \u00e8 divisibili in\" -> \" ed \u00e8 divisibile in\"\nconst divisibili = body.search(\", e teoricamente \u00e8 divisibili in\", { matchCase: true, matchWholeWord: false });\ndivisibili.load(\"items\");\nawait context.sync();\nif (divisibili.items.length > 0) {\n  divisibili.items[0].insertText(\" ed \u00e8 divisibile in\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 3) \", cio\u00e8 piccolissimo\" -> \", cio\u00e8 un piccolissimo\"\nconst cioe = body.search(\", cio\u00e8 piccolissimo\", { matchCase: true, matchWholeWord: false });\ncioe.load(\"items\");\nawait context.sync();\nif (cioe.items.length > 0) {\n  cioe.items[0].insertText(\", cio\u00e8 un piccolissimo\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 4) Move the document's \"_GoBack\" (last edit position) bookmark from its old\n//    spot to right after \"\u00e8 div\" inside the newly retyped \"divisibile\".\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst divSpot = body.search(\"\u00e8 div\", { matchCase: true, matchWholeWord: false });\ndivSpot.load(\"items\");\nawait context.sync();\nif (divSpot.items.length > 0) {\n  const collapsed = divSpot.items[0].getRange(Word.RangeLocation.end);\n  collapsed.insertBookmark(\"_GoBack\");\n}\nawait context.sync();\n", "ps1": "# \"modifica testo teoria avanzata\"\n#\n# Three textual edits on the document body:\n#  1) Paragraph \"Il tono e...\": the clause\n#       \", e teoricamente e divisibili in\"\n#     becomes\n#       \" ed e divisibile in\"\n#  2) Same paragraph, later: \", cioe piccolissimo\" becomes \", cioe un piccolissimo\"\n#  3) Paragraph \"Il semitono diatonico...\": the split run \"qua\"/\"ttro comma\"\n#     (with the _GoBack bookmark sitting between them) is normalized into a\n#     single \"quattro comma\" run.\n#\n# Doing (3) first lets Find/Replace naturally retire the pre-existing\n# _GoBack bookmark that lived inside \"qua|ttro comma\". We then re-create it\n# at its new resting place (inside \"divisibile\", right after \"e div\") to\n# match where Word leaves its \"last edit\" marker after the retype in\n# paragraph 2 -- Bookmarks.Add with the already-used name \"_GoBack\" simply\n# moves the single, uniquely named bookmark there.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# 1) Normalize \"qua\" + \"ttro comma\" into one run (text is unchanged).\nReplace-Text \"quattro comma\" \"quattro comma\"\n\n# 2) \", e teoricamente \u00e8 divisibili in\" -> \" ed \u00e8 divisibile in\"\nReplace-Text \", e teoricamente \u00e8 divisibili in\" \" ed \u00e8 divisibile in\"\n\n# 3) \", cio\u00e8 piccolissimo\" -> \", cio\u00e8 un piccolissimo\"\nReplace-Text \", cio\u00e8 piccolissimo\" \", cio\u00e8 un piccolissimo\"\n\n# 4) Move the document's \"_GoBack\" (last edit position) bookmark from its old\n#    spot to right after \"\u00e8 div\" inside the newly retyped \"divisibile\".\n$r = $d.Content.Duplicate\n$r.Find.ClearFormatting()\n$r.Find.Text = \"\u00e8 div\"\n$r.Find.Execute() | Out-Null\n$r.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $r) | Out-Null\n"}
